$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.654.61'
$ws.Range('E2').Value = '  -1.95%  '
$ws.Range('D3').Value = '1.792.72'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.01'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4451'
$ws.Range('E7').Value = '  +5.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3657'
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07278'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8518'
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.54'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.964.24'
$ws.Range('E12').Value = '  +7.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.604'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07064'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.82'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.261'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008645'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.78'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').Value = '26.695.26'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.135'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.76'
$ws.Range('E23').Value = '  -1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.981'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.74'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.38'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.169'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.160'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.21'
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08781'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.7394'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.151'
$ws.Range('E32').Value = '  -3.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.926'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.431'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01954'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05164'
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5261'
$ws.Range('E39').Value = '  +4.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.832'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.003'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1673'
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5080'
$ws.Range('E43').Value = '  +7.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.376'
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.43'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.951'
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.29'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.653'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06289'
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9109'
$ws.Range('E51').Value = '  -0.08%  '
